# Updated ppt and some things from module 3
#
# 1. Refresh the cached text of every auto-updating "date" placeholder
#    (the datetimeFigureOut field) on the Slide Master and on each of its
#    Custom Layouts, from 01/07/2020 -> 01/09/2020.
# 2. Nudge "Content Placeholder 3" on slide 9 down slightly
#    (y: 3886246 EMU -> 3891383 EMU).

$p = $ppt.ActivePresentation

$newDate = "01/09/2020"
$ppPlaceholderDate = 16
$EMU_PER_POINT = 12700

# --- 1. Slide Master date placeholder -------------------------------------
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    $isDatePh = $false
    try { $isDatePh = $sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate } catch { $isDatePh = $false }
    if ($isDatePh) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- 1b. Every Custom Layout's date placeholder ----------------------------
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        $isDatePh = $false
        try { $isDatePh = $sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate } catch { $isDatePh = $false }
        if ($isDatePh) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Move "Content Placeholder 3" on slide 9 -----------------------------
$s = $p.Slides.Item(9)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $leftEmu = [Math]::Round($sh.Left * $EMU_PER_POINT)
    $topEmu = [Math]::Round($sh.Top * $EMU_PER_POINT)
    if ($sh.Name -eq "Content Placeholder 3" -and $leftEmu -eq 6015519 -and $topEmu -eq 3886246) {
        $sh.Top = 3891383 / $EMU_PER_POINT
    }
}

Write-Output "edit applied"
